$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 703 ("「私できるよ」..." entry), shifting all rows below it up by one.
$ws.Rows.Item(703).Delete()
